$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update fecha2 value (row 3, col C): 30_05_2023 -> 21_07_2023 ---
$ws.Cells.Item(3, 3).Value = "21_07_2023"

# --- Row 5: finish "Conclusion Final" row by adding the C5 paragraph ---
$ws.Cells.Item(5, 3).Value = "Primera oracion de la conclusion final`nSeguna oracion de la conclusion final"
$ws.Cells.Item(5, 3).WrapText = $true
$ws.Rows.Item(5).RowHeight = 220.1

# --- Row 6: Introduccion ---
$ws.Cells.Item(6, 1).Value = "Introduccion"
$ws.Cells.Item(6, 2).Value = "Parrafo. Punto y aparte usando Ctrl+Enter"
$ws.Cells.Item(6, 3).Value = "Primer Renglon de la intro`nSegundo Renglon de la intro"
$ws.Cells.Item(6, 3).WrapText = $true
$ws.Rows.Item(6).RowHeight = 28.35

# --- Row 7: Objetivo ---
$ws.Cells.Item(7, 1).Value = "Objetivo"
$ws.Cells.Item(7, 2).Value = "Parrafo. Punto y aparte usando Ctrl+Enter"
$ws.Cells.Item(7, 3).Value = "Unica oracion del objetivo"
$ws.Rows.Item(7).RowHeight = 13.8

# --- Row 8: Metodologia ---
$ws.Cells.Item(8, 1).Value = "Metodologia"
$ws.Cells.Item(8, 2).Value = "Parrafo. Punto y aparte usando Ctrl+Enter"
$ws.Cells.Item(8, 3).Value = "Primer rengon de los objetivos`nSegundo renglon de los objetivos"
$ws.Cells.Item(8, 3).WrapText = $true
$ws.Rows.Item(8).RowHeight = 28.35

# --- Selection moves to C4 ---
[void]$ws.Range("C4").Select()
